$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores every "Price" cell in column D as plain text
# (t="inlineStr"), even when the text looks like a decimal number. Force the
# refreshed price cells whose new value would otherwise be auto-detected as a
# number by Excel to keep a text format, so they stay text cells like the rest
# of the column.
$textPriceCells = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D41", "D42", "D46", "D47", "D49"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Refreshed coin prices / 1h volume-change figures scraped this run, plus the
# PancakeSwap / InternetComputer(DFINITY) row-order swap (rows 24-25).
$ws.Range("D2").Value = '43.778.05'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '2.272.91'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '121.79'
$ws.Range("E5").Value = '  +7.65%  '
$ws.Range("D6").Value = '267.09'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = '0.651'
$ws.Range("E7").Value = '  +5.20%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = '0.627'
$ws.Range("E9").Value = '  +5.20%  '
$ws.Range("D10").Value = '48.60'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").Value = '0.0949'
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '9.28'
$ws.Range("E12").Value = '  +6.29%  '
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '15.65'
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").Value = '0.914'
$ws.Range("E15").Value = '  +7.10%  '
$ws.Range("D16").Value = '2.617.24'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").Value = '2.277.08'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '43.671.72'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("E19").Value = '  +3.15%  '
$ws.Range("D20").Value = '6.98'
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("D21").Value = '72.32'
$ws.Range("E21").Value = '  +1.58%  '
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '235.79'
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '2.91'
$ws.Range("E24").Value = '  +2.49%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '9.55'
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("D26").Value = '11.96'
$ws.Range("E26").Value = '  +5.74%  '
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D28").Value = '43.33'
$ws.Range("E28").Value = '  +6.46%  '
$ws.Range("D29").Value = '3.43'
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").Value = '174.07'
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("D32").Value = '21.71'
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").Value = '0.0929'
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("E34").Value = '  +2.85%  '
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").Value = '  +4.13%  '
$ws.Range("D36").Value = '4.27'
$ws.Range("E36").Value = '  +10.94%  '
$ws.Range("E37").Value = '  +10.08%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("E40").Value = '  +5.03%  '
$ws.Range("D41").Value = '74.05'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").Value = '13.76'
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("D46").Value = '5.89'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("D47").Value = '75.11'
$ws.Range("E47").Value = '  +43.32%  '
$ws.Range("E48").Value = '  +4.04%  '
$ws.Range("D49").Value = '103.45'
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("E51").Value = '  -1.58%  '
